$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 269, shifting rows 269:347 down to 270:348
$ws.Rows.Item(269).Insert()

# Populate the newly inserted row 269 with its values
$ws.Cells.Item(269, 1).Value = 10
$ws.Cells.Item(269, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(269, 3).Value = "La Araucanía"
$ws.Cells.Item(269, 4).Value = 44988
$ws.Cells.Item(269, 5).Value = 9
$ws.Cells.Item(269, 6).Value = 100112039
$ws.Cells.Item(269, 7).Value = "Ciboulette"
$ws.Cells.Item(269, 8).Value = "Sin especificar"
$ws.Cells.Item(269, 9).Value = "Primera"
$ws.Cells.Item(269, 10).Value = 30
$ws.Cells.Item(269, 11).Value = 5000
$ws.Cells.Item(269, 12).Value = 5000
$ws.Cells.Item(269, 13).Value = 5000
$ws.Cells.Item(269, 14).Value = "$/docena de atados"
$ws.Cells.Item(269, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(269, 16).Value = 1667
$ws.Cells.Item(269, 17).Value = 3
$ws.Cells.Item(269, 18).Value = "Hortaliza"
